$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BundleID value in B2
$ws.Range("B2").Value = "EB000016"

# Add new OpportunityID column
$ws.Range("C1").Value = "OpportunityID"
$ws.Range("C2").Value = "OPE-0002907630"

# Update selection to C5 as seen in target file
$ws.Range("C5").Select()
